$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Index")

# Update B9: append " model results" to the existing description text
$ws.Range("B9").Value = "Code for manuscript- only inlcudes last survey wave model results"

# Fill in row 10 with the new file path / description entries
$ws.Range("A10").Value = "/Users/carmenrodriguez/Desktop/Research Projects/BayesBinMix/ecbayesbinmix/ECbayesbinmix_manuscript.R"
$ws.Range("B10").Value = "Code for manuscript- EC data and Bayesbinmix resutls"

# B10 picks up the wrap-text bordered style used elsewhere in the table (same as column B cells above)
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row height matches the other wrapped rows
$ws.Rows.Item(10).RowHeight = 34

# Move the active selection to C9 (matches saved selection state)
$ws.Range("C9").Select()
